$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("D4:D9")
$fcs = $rng.FormatConditions
$fc1 = $fcs.Item(1)
$fc1.Font.Color = 255
Write-Output "done"
